$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings are stored as text
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.232.73"
$ws.Range("E2").Value = "  -2.30%  "

$ws.Range("D3").Value = "1.876.14"
$ws.Range("E3").Value = "  -2.18%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "235.42"
$ws.Range("E5").Value = "  -1.56%  "

$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("D7").Value = "0.4829"
$ws.Range("E7").Value = "  -1.85%  "

$ws.Range("D8").Value = "0.2867"
$ws.Range("E8").Value = "  -3.29%  "

$ws.Range("D9").Value = "0.06572"
$ws.Range("E9").Value = "  -2.97%  "

$ws.Range("D10").Value = "1.888.85"
$ws.Range("E10").Value = "  -1.56%  "

$ws.Range("E11").Value = "  -2.06%  "

$ws.Range("D12").Value = "0.07322"
$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("D13").Value = "5.126"
$ws.Range("E13").Value = "  -0.37%  "

$ws.Range("D14").Value = "86.88"
$ws.Range("E14").Value = "  -3.52%  "

$ws.Range("D15").Value = "0.6515"
$ws.Range("E15").Value = "  -2.97%  "

$ws.Range("D16").Value = "30.208.75"
$ws.Range("E16").Value = "  -2.23%  "

$ws.Range("D17").Value = "13.29"
$ws.Range("E17").Value = "  -0.91%  "

$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D19").Value = "0.000007721"
$ws.Range("E19").Value = "  -2.82%  "

$ws.Range("D20").Value = "2.120.12"
$ws.Range("E20").Value = "  -1.49%  "

$ws.Range("D21").Value = "5.342"
$ws.Range("E21").Value = "  +3.00%  "

$ws.Range("D22").Value = "0.9981"
$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("D23").Value = "192.58"
$ws.Range("E23").Value = "  -6.56%  "

$ws.Range("D24").Value = "6.090"
$ws.Range("E24").Value = "  -2.27%  "

$ws.Range("D25").Value = "9.240"
$ws.Range("E25").Value = "  -4.61%  "

$ws.Range("D26").Value = "161.99"
$ws.Range("E26").Value = "  +2.65%  "

$ws.Range("E27").Value = "  -4.96%  "

$ws.Range("D28").Value = "1.904"
$ws.Range("E28").Value = "  -3.67%  "

$ws.Range("D29").Value = "1.433"
$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("E30").Value = "  -1.68%  "

$ws.Range("D31").Value = "0.09074"
$ws.Range("E31").Value = "  -1.08%  "

$ws.Range("D32").Value = "3.994"
$ws.Range("E32").Value = "  -1.76%  "

$ws.Range("D33").Value = "0.05058"
$ws.Range("E33").Value = "  -2.32%  "

$ws.Range("D34").Value = "0.7087"
$ws.Range("E34").Value = "  -5.75%  "

$ws.Range("E35").Value = "  -2.57%  "

$ws.Range("D36").Value = "2.700"
$ws.Range("E36").Value = "  -1.38%  "

$ws.Range("D37").Value = "0.01778"
$ws.Range("E37").Value = "  -3.89%  "

$ws.Range("D38").Value = "2.637"
$ws.Range("E38").Value = "  -3.69%  "

$ws.Range("D39").Value = "0.9196"
$ws.Range("E39").Value = "  -0.48%  "

$ws.Range("D40").Value = "2.032"
$ws.Range("E40").Value = "  -3.08%  "

$ws.Range("D41").Value = "105.65"
$ws.Range("E41").Value = "  -1.46%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.4249"
$ws.Range("E42").Value = "  -6.00%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.773"
$ws.Range("E43").Value = "  -1.89%  "

$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.84%  "

$ws.Range("D45").Value = "7.357"
$ws.Range("E45").Value = "  -4.51%  "

$ws.Range("D46").Value = "0.1309"
$ws.Range("E46").Value = "  -6.51%  "

$ws.Range("D47").Value = "64.67"
$ws.Range("E47").Value = "  -3.14%  "

$ws.Range("D48").Value = "8.905"
$ws.Range("E48").Value = "  -0.40%  "

$ws.Range("D49").Value = "0.05743"
$ws.Range("E49").Value = "  -3.54%  "

$ws.Range("D50").Value = "33.58"
$ws.Range("E50").Value = "  -4.26%  "

$ws.Range("D51").Value = "0.3800"
$ws.Range("E51").Value = "  -7.09%  "

# Restore normal style/number format on column D
$ws.Range("D2:D51").Style = "Normal"
